$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-19) of this sheet get reordered: the new row N contains
# the data that used to live in row Map[N] (a permutation of 2..19).
$map = @{2=3; 3=18; 4=14; 5=17; 6=10; 7=16; 8=8; 9=11; 10=7; 11=19; 12=5; 13=2; 14=6; 15=9; 16=12; 17=13; 18=4; 19=15}

$firstRow = 2
$lastRow = 19
$lastCol = 20  # column T

# 1) Snapshot every cell of every data row (by Value2, which round-trips
#    numbers/dates/strings faithfully) before we overwrite anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each new row from the snapshot of its source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
